# Record that LeetCode problem 413 (a DP problem) was solved: append a new
# entry (row 23) to the "新题" (new problems) sheet with today's date, the
# problem name, and a "done" status - mirroring the existing rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# Copy the date cell's formatting from the row above so the new date cell
# keeps the same (already-existing) date style instead of minting a new one.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A23").Value = 43551
$ws.Range("B23").Value = "413 dp"
$ws.Range("F23").Value = "done"

$ws.Range("A23").Select()
